$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1, matching the style of the existing header row (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-23
$data = @(
    @(7, 8),
    @(6, 8),
    @(6, 6),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(6, 7),
    @(5, 6),
    @(11, 11),
    @(6, 6),
    @(8, 8),
    @(5, 6),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(6, 8),
    @(6, 8),
    @(6, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
